# Update the cryptos worksheet with refreshed price / 1h-volume-change figures
# (and swap the Stacks/Filecoin row order) per the Fri Jun 14 05:35:48 UTC 2024
# GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.857.31"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "3.510.34"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'602.91"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("D6").Value = "'147.64"
$ws.Range("E6").Value = "  -3.07%  "
$ws.Range("D7").Value = "3.509.92"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("D12").Value = "'0.423"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").Value = "'0.0000215"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").Value = "4.102.43"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "'31.55"
$ws.Range("E15").Value = "  -3.26%  "
$ws.Range("D16").Value = "3.503.29"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "66.865.98"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").Value = "'10.65"
$ws.Range("E19").Value = "  +7.56%  "
$ws.Range("D20").Value = "'6.40"
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("D21").Value = "'15.40"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").Value = "'434.92"
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("E23").Value = "  -3.44%  "
$ws.Range("D24").Value = "'79.87"
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("D25").Value = "3.643.16"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").Value = "'0.0000120"
$ws.Range("E27").Value = "  -5.01%  "
$ws.Range("D28").Value = "'9.85"
$ws.Range("E28").Value = "  -2.51%  "
$ws.Range("D29").Value = "'8.29"
$ws.Range("E29").Value = "  -6.00%  "
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("E31").Value = "  -3.91%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("D34").Value = "'25.32"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("D35").Value = "3.499.57"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  -4.05%  "
$ws.Range("D37").Value = "'5.89"
$ws.Range("E37").Value = "  -4.79%  "
$ws.Range("D38").Value = "'8.03"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "'0.0891"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "'169.12"
$ws.Range("E42").Value = "  -2.49%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'5.43"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'2.09"
$ws.Range("E44").Value = "  -9.73%  "
$ws.Range("D45").Value = "'0.898"
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("D46").Value = "'29.05"
$ws.Range("E46").Value = "  -3.65%  "
$ws.Range("D47").Value = "'45.72"
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Value = "'7.47"
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("D50").Value = "'2.42"
$ws.Range("E50").Value = "  -4.35%  "
$ws.Range("D51").Value = "'0.983"
$ws.Range("E51").Value = "  -1.17%  "
